$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row for the beverage with Id 53fad91e-c4c8-42f5-a81f-f809f838c37f
# (simulates clicking the "Delete" button for that row), which shifts the
# remaining rows up.
$ws.Rows.Item(2).Delete()

# Add the newly created beverage entry as a new row at the end of the data.
$newRow = 3

# Force the Date/Time columns to be stored as plain text (matching the
# existing rows), rather than letting Excel auto-convert them to date/time
# serial numbers.
$ws.Range("G" + $newRow + ":H" + $newRow).NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "32b13210-7cf0-4040-8158-3648e246efed"
$ws.Cells.Item($newRow, 2).Value = ",s,xkks"
$ws.Cells.Item($newRow, 3).Value = "Water"
$ws.Cells.Item($newRow, 4).Value = "7UP"
$ws.Cells.Item($newRow, 5).Value = 10
$ws.Cells.Item($newRow, 6).Value = 100
$ws.Cells.Item($newRow, 7).Value = "2024-09-14"
$ws.Cells.Item($newRow, 8).Value = "12:32:28"
